# Target state (row 1): A1="s", B1="Lâmpada", C1=0, D1=FALSE, F1=FALSE (E1 left empty)
# and the sheet's used range/dimension grows from A1:A1 to A1:F1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep cells unstyled (no explicit style index), matching the source sheet's
# plain cells, by resetting to the "Normal" style before/while writing values.
$ws.Range("A1").Style = "Normal"
$ws.Range("B1").Style = "Normal"
$ws.Range("C1").Style = "Normal"
$ws.Range("D1").Style = "Normal"
$ws.Range("F1").Style = "Normal"

$ws.Range("A1").Value = "s"
$ws.Range("B1").Value = "Lâmpada"
$ws.Range("C1").Value = 0
$ws.Range("D1").Value = $false
$ws.Range("F1").Value = $false
